$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 260; this shifts the existing rows 260-346
# down to 261-347 (carrying their values/styles with them), and extends
# the sheet dimension to A1:R347 automatically.
$ws.Rows(260).Insert()

# Populate the freshly inserted row 260 with the new data record. The
# descriptive columns (market/region/category/etc.) match the record
# that used to sit at row 260 (now at row 261); only the date and the
# measurement columns carry new values.
$ws.Range("A260").Value = 3
$ws.Range("B260").Value = "Femacal de La Calera"
$ws.Range("C260").Value = "Coquimbo"
$ws.Range("D260").Value = 44627
$ws.Range("E260").Value = 5
$ws.Range("F260").Value = 100112017
$ws.Range("G260").Value = "Apio"
$ws.Range("H260").Value = "Americana (o)"
$ws.Range("I260").Value = "Primera"
$ws.Range("J260").Value = 282
$ws.Range("K260").Value = 9500
$ws.Range("L260").Value = 10500
$ws.Range("M260").Value = 9995
$ws.Range("N260").Value = "$/docena de matas"
$ws.Range("O260").Value = "Pan de Azúcar"
$ws.Range("P260").Value = 1666
$ws.Range("Q260").Value = 6
$ws.Range("R260").Value = "Hortaliza"
